$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17-27 down to 18-28
$ws.Rows(17).Insert()

# Populate the new row 17 with the same categorical data as the surrounding rows,
# plus the new weekly observation values.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 45236
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100104
$ws.Cells.Item(17, 8).Value = "Frutos de pepita"
$ws.Cells.Item(17, 9).Value = 100104004
$ws.Cells.Item(17, 10).Value = "Níspero"
$ws.Cells.Item(17, 11).Value = "Californiana(o)"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 30
$ws.Cells.Item(17, 14).Value = 30000
$ws.Cells.Item(17, 15).Value = 30000
$ws.Cells.Item(17, 16).Value = 30000
$ws.Cells.Item(17, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 19).Value = 3000
$ws.Cells.Item(17, 20).Value = 10
